$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new daily data row (row 40) with date serial + counts
$ws.Range("A40").Value = 45995
$ws.Range("B40").Value = 666
$ws.Range("C40").Value = 8
$ws.Range("D40").Value = 658

# Update selection to reflect the newly entered row, matching the saved view state
$ws.Range("A40:D40").Select()
